# Translate the "Review / Date / Location" paragraphs of the document to
# Spanish (es-AR), matching the author's commit: each paragraph (and every
# run within it) is tagged with w:lang="es-AR" via the paragraph mark's
# rPr and each run's rPr. We rebuild each paragraph's contents precisely
# (including the paragraph-mark formatting, which plain Range.Text /
# Range.LanguageID assignment cannot reach) using Range.InsertXML, which
# replaces a range's contents with arbitrary WordprocessingML - exactly
# like pasting flattened OOXML into that range.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Paragraph 2: "Review: ..." -> "Reseña: ..." -----------------------
# Two runs: the translated review sentence, then a separate run holding
# just the trailing period (mirrors the target markup's run split).
$reviewXml = '<w:p>' +
    '<w:pPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr>' +
    '<w:t>Reseña: A menudo elijo Fourth Coffee como lugar para mis reuniones con clientes por las mañanas entre semana. Tengo un pequeño negocio, y las personas que trabajan en Fourth Coffee siempre son muy amables. Dejan una buena impresión en mis clientes. Además, hay muchas opciones de bebidas, buen wi-fi y suficientes asientos. Algunos de mis cafés favoritos son el latte de miel con lavanda y, en invierno, el latte de manzana y chai. También ofrecen deliciosos productos horneados</w:t>' +
    '</w:r>' +
    '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>.</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs(2).Range.InsertXML($pkgHeader + $reviewXml + $pkgFooter)

# --- Paragraph 3: "Date: ..." -> "Fecha: ..." ---------------------------
$dateXml = '<w:p>' +
    '<w:pPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Fecha: 21 de octubre de 2018</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs(3).Range.InsertXML($pkgHeader + $dateXml + $pkgFooter)

# --- Paragraph 4: "Location: ..." -> "Ubicación: ..." -------------------
$locationXml = '<w:p>' +
    '<w:pPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Ubicación: Chicago, Illinois</w:t></w:r>' +
    '</w:p>'
$d.Paragraphs(4).Range.InsertXML($pkgHeader + $locationXml + $pkgFooter)

# --- Remove the trailing empty paragraph (paragraph 5) ------------------
# Extend the delete range back to swallow the preceding paragraph mark so
# the empty paragraph collapses away entirely (deleting just its own,
# empty range leaves a stray paragraph mark behind).
$p5 = $d.Paragraphs(5)
$delRange = $d.Range($p5.Range.Start - 1, $p5.Range.End)
$delRange.Delete()
